$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (V1:Y1)
$ws.Range("V1").Value = "ATTA H Cost for Failures"
$ws.Range("W1").Value = "ATTA R Cost for Failures"
$ws.Range("X1").Value = "Tsarouchi MIN H Cost for Failures"
$ws.Range("Y1").Value = "Tsarouchi MIN R Cost for Failures"

# New data values for rows 2-11
$data = @(
    @(2.6720876217459, 12.3283111690625, 12.783333333333299, 5.3020987369495503),
    @(0.74408058723098702, 18.3688860622863, 13.8666666666666, 5.6549112458010402),
    @(2.6865228560818499, 13.3286065976421, 10.8333333333333, 6.0163391306524998),
    @(2.9310042475495801, 13.579284661867799, 18.2823849248438, 6.0088646189847603),
    @(1.7525117617388699, 15.182411123623201, 14.3, 5.9413386770403402),
    @(1.05604380662411, 16.093658535388499, 13.216666666666599, 5.9345930505849802),
    @(2.55449287289477, 14.920576240820999, 14.733333333333301, 6.5229384990264299),
    @(3.5554878088780999, 14.868615354645099, 14.3, 6.3621642797684199),
    @(1.3740478201514501, 13.595029701367601, 13, 4.5889100553602802),
    @(2.04408086108378, 15.5826476926888, 14.733333333333301, 5.8068940971754204)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 22).Value = $data[$i][0]
    $ws.Cells.Item($row, 23).Value = $data[$i][1]
    $ws.Cells.Item($row, 24).Value = $data[$i][2]
    $ws.Cells.Item($row, 25).Value = $data[$i][3]
}

# Update selection to match target state
$ws.Range("N18").Select()
